$wb = $excel.ActiveWorkbook

# zh-cn sheet: update Correspond Handoff Datetime (D2) and Correspond Handback DateTime (G2)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D2").Value = "2016-02-17 09:09:36"
$wsZh.Range("G2").Value = "2016-02-17 09:10:24"

# de-de sheet: update Correspond Handoff Datetime (D2) and Correspond Handback DateTime (G2)
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D2").Value = "2016-02-17 09:09:47"
$wsDe.Range("G2").Value = "2016-02-17 09:10:51"
